$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha")
$ws.Activate()

# --- M2: replace the requirement write-up with the new text ---
$m2Text = @'
🧾 Requisito Definido
O sistema IRIS deve capturar o campo de status que já é disponibilizado pela API do software NxTDC durante a chamada de interface. Este status deve ser usado para formatar o nome final do arquivo jobReport.xml antes de disponibilizá-lo na pasta para o sistema AIT.
•	Formato do Novo Nome: <nome>_<status>.xml
•	Escopo: A alteração deve ocorrer somente na interface do IRIS com a TDC.
•	Fora do Escopo: Nenhuma alteração deve ser feita na parte do SAP.
🧑‍💼 Contato
•	Product Owner (Negócio): Izabel (Contato para dúvidas de negócio e priorização).
•	Referência Técnica (Interface): Pablo (Contexto técnico sobre a interface NxTDC que ele mexeu).
📊 Priorização
•	Nível: Média.
•	Justificativa (Valor de Negócio): Esta é uma solicitação direta da área de negócio que visa melhorar a rastreabilidade e a eficiência operacional do sistema AIT. Ao incluir o status no nome do arquivo, o AIT poderá automatizar o processamento subsequente (ex: mover arquivos "Fail" para uma pasta de erro) sem a necessidade de abrir e processar o conteúdo de cada XML para descobrir o resultado. Isso reduz o tempo de processamento e simplifica o monitoramento de falhas.
🔍 Processo Atual (AS IS)
1.	O IRIS chama a interface do NxTDC para buscar o arquivo jobReport.xml.
2.	A API do NxTDC retorna o status do job, mas o IRIS não utiliza (ignora) essa informação.
3.	O IRIS disponibiliza o arquivo com o nome estático jobReport.xml em uma pasta.
4.	O AIT localiza o arquivo jobReport.xml e disponibiliza o arquivo para acesso do usuário final..
🚀 Processo Futuro (TO BE)
1.	O IRIS chama a interface do NxTDC.
2.	O IRIS captura o valor do status retornado pela API.
3.	O IRIS usa esse status para formatar o nome do arquivo de relatório (ex: JobReport_Success.xml, JobReport_Fail.xml, ou JobReport.xml se o status vier em branco).
4.	O IRIS disponibiliza o arquivo já renomeado na pasta.
5.	O AIT localiza o arquivo e identifica o status imediatamente pelo nome, direcionando seu fluxo de trabalho.
🎯 Objetivo da Mudança
Permitir que o sistema consumidor (AIT) identifique o resultado de um job (Sucesso ou Falha) pela nomenclatura do arquivo, otimizando a automação e o monitoramento do processo sem a necessidade de ler o conteúdo do jobReport.xml.
✅ Critérios de Aceite (BDD)
Cenário 1: Job com Sucesso
•	Dado que o IRIS chamou a interface do NxTDC
•	E a API do NxTDC retornou um status de "Success" (ou equivalente a Sucesso)
•	Quando o IRIS for disponibilizar o arquivo de relatório na pasta
•	Então o arquivo deve ser nomeado como JobReport_Success.xml.
Cenário 2: Job com Falha
•	Dado que o IRIS chamou a interface do NxTDC
•	E a API do NxTDC retornou um status de "Fail" (ou equivalente a Erro)
•	Quando o IRIS for disponibilizar o arquivo de relatório na pasta
•	Então o arquivo deve ser nomeado como JobReport_Fail.xml.
Cenário 3: Status em Branco (Legado ou Nulo)
•	Dado que o IRIS chamou a interface do NxTDC
•	E a API do NxTDC retornou um status em branco ou nulo
•	Quando o IRIS for disponibilizar o arquivo de relatório na pasta
•	Então o arquivo deve ser nomeado como JobReport.xml (sem sufixo, conforme sugestão de "Status em branco").
Cenário 4: Integridade do SAP
•	Dado que esta funcionalidade foi implementada no IRIS
•	Quando qualquer processo SAP relacionado for executado
•	Então o comportamento do SAP deve permanecer inalterado (conforme escopo definido).
🧪 Cenários de Teste de Validação (Caminho Feliz)
1.	Validação de Sucesso (Status "Success"):
o	Passos: Simular (via mock ou teste) uma chamada à API do NxTDC onde o status retornado é "Success". Executar o fluxo do IRIS.
o	Esperado: Verificar na pasta de destino que o arquivo jobReport.xml foi salvo com o nome JobReport_Success.xml.
2.	Validação de Falha (Status "Fail"):
o	Passos: Simular uma chamada à API do NxTDC onde o status retornado é "Fail". Executar o fluxo do IRIS.
o	Esperado: Verificar na pasta de destino que o arquivo jobReport.xml foi salvo com o nome JobReport_Fail.xml.
❌ Cenários de Teste de Rejeição (Exceções)
1.	Rejeição de Status (Nulo ou Vazio):
o	Passos: Simular uma chamada à API do NxTDC onde o campo de status retorna null ou uma string vazia "".
o	Esperado: Verificar na pasta de destino que o arquivo foi salvo com o nome legado JobReport.xml (sem sufixo _), conforme Critério de Aceite 3.
2.	Rejeição de Status (Não Mapeado):
o	Passos: Simular uma chamada à API do NxTDC onde o status retorna um valor inesperado (ex: "Warning", "Processing", "Cancelled").
o	Esperado: (Ponto de Atenção para Refinamento) O requisito não define esse comportamento. A US deve ser atualizada após alinhamento com a PO (Izabel) e o time. 
	(Sugestão de regra): Se o status não for "Success" ou "Fail", tratar como "Fail" (ex: JobReport_Fail.xml) ou como um status desconhecido (ex: JobReport_Unknown.xml) para garantir que não seja processado como sucesso.
3.	Teste de Regressão (SAP):
o	Passos: Executar o fluxo de ponta a ponta, incluindo os processos SAP que tangenciam essa interface.
o	Esperado: Confirmar que não houve nenhuma quebra ou alteração de comportamento no lado do SAP, conforme requisito.
4.	Teste de Dependência (Impacto no AIT):
o	Passos: Executar o fluxo de Sucesso (gerando JobReport_Success.xml). Verificar o comportamento do sistema AIT (consumidor).
o	Esperado: (Alerta de Dependência para o PO) O AIT precisa estar preparado para ler os novos nomes de arquivo. Se o AIT ainda estiver procurando exclusivamente por jobReport.xml, o fluxo dele irá falhar. Esta entrega precisa ser coordenada com o time do AIT.
'@
$ws.Range("M2").Value = $m2Text

# --- Updated Jira ticket references ---
$ws.Range("R2").Value = "TRE-1626"
$ws.Range("S2").Value = "TRE-1627"
$ws.Range("T2").Value = "TRE-1628"

# --- View state: selection moved to S9, sheet scrolled so column L is leftmost ---
$ws.Range("L1").Select()
$excel.ActiveWindow.ScrollColumn = 12
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("S9").Select()
